$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 ---
$ws.Range("M2").Value = 2.740639666666667
$ws.Range("N2").Value = 8.221919
$ws.Range("O2").Value = 0.05814502584416985
$ws.Range("P2").Value = 0.05814502584416984
$ws.Range("Q2").Value = 0.1864091746611111
$ws.Range("R2").Value = 1.67768257195
$ws.Range("S2").Value = 0.05814502584416985
$ws.Range("T2").Value = 0.05814502584416984

# --- Update row 3 ---
$ws.Range("O3").Value = 0.7843079965148284
$ws.Range("P3").Value = 0.7843079965148283
$ws.Range("S3").Value = 0.7843079965148284
$ws.Range("T3").Value = 0.7843079965148283

# --- Update row 4 ---
$ws.Range("M4").Value = 7.396246333333333
$ws.Range("N4").Value = 22.188739
$ws.Range("O4").Value = 0.1569177223230415
$ws.Range("P4").Value = 0.1569177223230415
$ws.Range("Q4").Value = 0.5030680214388888
$ws.Range("R4").Value = 4.52761219295
$ws.Range("S4").Value = 0.1569177223230415
$ws.Range("T4").Value = 0.1569177223230415

# --- Add new row 5 ---
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Ihh"
$ws.Range("C5").Value = "Boc"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.06801666666666667
$ws.Range("H5").Value = 0.20405
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.02965966666666667
$ws.Range("N5").Value = 0.088979
$ws.Range("O5").Value = 0.0006292553179602461
$ws.Range("P5").Value = 0.0006292553179602461
$ws.Range("Q5").Value = 0.002017351661111111
$ws.Range("R5").Value = 0.01815616495
$ws.Range("S5").Value = 0.0006292553179602461
$ws.Range("T5").Value = 0.0006292553179602461
